$wb = $excel.ActiveWorkbook

$wsPeople = $wb.Worksheets.Item("People")
$wsChildren = $wb.Worksheets.Item("Children")

# --- Sheet "People": insert a new header/label row above the table header ---
$wsPeople.Rows.Item(2).Insert()
$wsPeople.Range("B2").Value = "Employee"
$wsPeople.Range("C2").Value = "Person"

# Update the active selection on the People sheet without changing the
# workbook's active sheet (Children stays active, as in the original file).
[void]$wsPeople.Activate()
[void]$wsPeople.Range("C2").Select()
[void]$wsChildren.Activate()

# --- Sheet "Children": add a new "Age" column (D), and fill in a
#     previously-missing Hobbies value for row 4 ---
$wsChildren.Range("D1").Value = "Age"
$wsChildren.Range("D2").Value = 1
$wsChildren.Range("D3").Value = 2
$wsChildren.Range("C4").Value = "Movies"
$wsChildren.Range("D5").Value = "."
$wsChildren.Range("D4").Value = "-"
$wsChildren.Range("D6").Value = 7
